# Updated cryptos list on Sun Jan  7 18:42:09 UTC 2024 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for each
# coin row (2-51) to the latest scraped values. D-column entries that look
# numeric (e.g. "306.30") are explicitly formatted as text first so Excel
# doesn't silently coerce them to numbers and drop significant trailing
# digits / thousands separators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.610.94"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "2.243.96"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.30"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.19"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.30"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "2.585.58"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "2.233.49"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.834"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "44.417.01"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.89"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.32"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.28"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +7.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.77"
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.26"
$ws.Range("E29").Value = "  -3.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.92"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.94"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.13"
$ws.Range("E32").Value = "  -1.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0788"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("E38").Value = "  +5.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.29"
$ws.Range("E39").Value = "  +4.55%  "
$ws.Range("E40").Value = "  -5.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.79"
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0300"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "1.813.15"
$ws.Range("E44").Value = "  +3.32%  "
$ws.Range("E45").Value = "  +12.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "81.72"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.188"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.49"
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.85"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.54"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.22"
$ws.Range("E51").Value = "  -1.27%  "
